# Fix bug: "Non-consistent bid amounts when clicking on landing page and
# add bid page". Updated on the "Bug Log" sheet of the Bug Metrics
# workbook under item 27 (row 29):
#   - Row 28, col C: fix the "Boostrap " typo -> "Bootstrap "
#   - Insert a new logged bug as row 29, copied/styled like row 28, for
#     the new "Non-consistent bid amounts..." entry (S/N 27, iteration 3,
#     Function "Bootstrap ", Status "Resolved", discovered/resolved on
#     14/11/2019, fixed by "Matthew & DaEun").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")
$ws.Activate()

# --- Fix the typo in the existing row 28 "Function" cell ---
$ws.Cells.Item(28, 3).Value = "Bootstrap "

# --- Build out new row 29, matching row 28's look & feel ---
# Columns A, C:H inherit their formatting from row 28. Column B has no
# thick-border formatting (it mirrors the sheet's blank filler rows, e.g.
# G19), so pull its format from there instead of from B28.
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)

$ws.Range("G19").Copy()
$ws.Range("B29").PasteSpecial(-4122)

$ws.Range("C28:H28").Copy()
$ws.Range("C29:H29").PasteSpecial(-4122)

# --- Populate the new row's data ---
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(29, 3).Value = "Bootstrap "
$ws.Cells.Item(29, 4).Value = "Non-consistent bid amounts when clicking on landing page and add bid page "
$ws.Cells.Item(29, 5).Value = "Resolved"
$ws.Cells.Item(29, 6).Value = "14/11/2019"
$ws.Cells.Item(29, 7).Value = "14/11/2019"
$ws.Cells.Item(29, 8).Value = "Matthew & DaEun"

# --- Match the author's final selection/scroll position on the sheet ---
$ws.Range("G29").Select()
